# Updates the "cryptos" price table with the latest scraped figures.
# Price cells (column D) are plain display text (e.g. "58.268.17",
# "1.00", "0.0000165") that Excel would otherwise reinterpret as numbers
# (stripping trailing zeros / mangling thousands-separated values), so
# each is forced to Text format before the assignment and the style is
# reset back to Normal afterwards to avoid leaving a number-format
# override on the cell. Volume(%) cells (column E) are plain text
# already (they contain spaces/`%`) and need no special handling. Rows
# 46/47 additionally had their Coin/Link values swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.268.17'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.140.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('E6').Value = '  +0.93%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.138.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.469'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('E12').Value = '  +4.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.679.28'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000165'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '58.384.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.142.33'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '362.01'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.20'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.166'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0879'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.37'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.22%  '
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.51'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.13'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.64%  '
$ws.Range('E34').Value = '  -2.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.91'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.11'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.93'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.50%  '
$ws.Range('E38').Value = '  +1.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.68'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0675'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.516.86'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.703'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.02'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '37.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.182.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0269'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.998'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.08%  '
$ws.Range('E50').Value = '  -1.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.745'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.49%  '
